$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.626.58'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '1.641.68'
$ws.Range("E3").Value = '  +0.64%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.39%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.02'
$ws.Range("E5").Value = '  +0.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.504'
$ws.Range("E6").Value = '  +0.95%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.01'
$ws.Range("E7").Value = '  +0.34%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.252'
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0626'
$ws.Range("E9").Value = '  +0.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.22'
$ws.Range("E10").Value = '  +0.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0841'
$ws.Range("E11").Value = '  +0.02%  '
$ws.Range("D12").Value = '1.871.83'
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("D13").Value = '1.680.35'
$ws.Range("E13").Value = '  +3.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.19'
$ws.Range("E14").Value = '  +2.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.529'
$ws.Range("E15").Value = '  +1.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.36'
$ws.Range("E16").Value = '  +3.10%  '
$ws.Range("D17").Value = '26.688.63'
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("D18").Value = '0.0₃0746'
$ws.Range("E18").Value = '  +0.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.71'
$ws.Range("E19").Value = '  -1.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.01'
$ws.Range("E20").Value = '  +0.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.36'
$ws.Range("E21").Value = '  +1.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.29'
$ws.Range("E22").Value = '  +2.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.51'
$ws.Range("E23").Value = '  +1.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.20'
$ws.Range("E24").Value = '  +12.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.75'
$ws.Range("E25").Value = '  -1.58%  '
$ws.Range("E26").Value = '  +0.50%  '
$ws.Range("E27").Value = '  -0.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.15'
$ws.Range("E28").Value = '  +4.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.73'
$ws.Range("E29").Value = '  +1.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0515'
$ws.Range("E30").Value = '  +2.10%  '
$ws.Range("E31").Value = '  +0.25%  '
$ws.Range("E32").Value = '  +2.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.04'
$ws.Range("E33").Value = '  +2.07%  '
$ws.Range("D34").Value = '1.273.71'
$ws.Range("E34").Value = '  +4.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.53'
$ws.Range("E35").Value = '  +2.17%  '
$ws.Range("E36").Value = '  +5.67%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.41'
$ws.Range("E37").Value = '  +0.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.530'
$ws.Range("E38").Value = '  +5.98%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.825'
$ws.Range("E39").Value = '  +2.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.01'
$ws.Range("E40").Value = '  +0.41%  '
$ws.Range("E41").Value = '  +2.76%  '
$ws.Range("E42").Value = '  -1.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.45'
$ws.Range("E43").Value = '  +2.32%  '
$ws.Range("D44").Value = '1.781.86'
$ws.Range("E44").Value = '  +0.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.28'
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '59.75'
$ws.Range("E46").Value = '  +8.53%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.60'
$ws.Range("E47").Value = '  +2.31%  '
$ws.Range("D48").Value = '0.0₆0104'
$ws.Range("E48").Value = '  -0.16%  '
$ws.Range("E49").Value = '  +0.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.81'
$ws.Range("E50").Value = '  +2.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0969'
$ws.Range("E51").Value = '  +3.00%  '
